$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: make the "Ejercicio 6" title and its column headers bold ---
# (row 2 first, then row 1 -- matches the style-table ordering seen in the
# saved workbook: a "bold only" xf followed by a "bold + centered" xf)
$ws.Range("A2:E2").Font.Bold = $true
$ws.Range("A1:E1").Font.Bold = $true

# --- Add "Ejercicio 5" block starting at row 9 ---
$ws.Cells.Item(9, 1).Value = "Ejercicio 5"

$ws.Cells.Item(10, 1).Value = "Nivel"
$ws.Cells.Item(10, 2).Value = "Valor"
$ws.Cells.Item(10, 3).Value = "Diferencia"

$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = 4

$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = 12

$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = 52

$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 2).Value = 252

# Column C (formula description) filled top-to-bottom after the numbers
$ws.Cells.Item(12, 3).Value = "f(n-1) +8"
$ws.Cells.Item(13, 3).Value = "f(n-1)+ 40"
$ws.Cells.Item(14, 3).Value = "f(n-1)+200"

# Column D (difference as a product) filled top-to-bottom last
$ws.Cells.Item(12, 4).Value = "8*1"
$ws.Cells.Item(13, 4).Value = "8*5"
$ws.Cells.Item(14, 4).Value = "8*25"

# --- Page setup tweak captured in the saved file ---
$ws.PageSetup.PaperSize = 10000
$ws.PageSetup.Orientation = 1

# --- Final selection left on the sheet ---
$null = $ws.Range("F15").Select()
